$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet's display name (sheet tab), keep internal data same.
$ws.Name = "BetaFiberA"

# Correct a floating point value in J15 (higher precision recomputation).
$ws.Range("J15").Value = 0.9986981690986969

# Add new row 16 with data for HKL index 14 / "HexGrid-60degTilt5degRes".
# Copy A15's formatting (bold/centered/bordered "index" style) down to A16,
# then overwrite the copied value with the new row's index.
$ws.Range("A15").Copy($ws.Range("A16"))
$ws.Range("A16").Value = 14

$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"

$ws.Range("C16").Value = 1.027139640605709
$ws.Range("D16").Value = 0.9445305089422378
$ws.Range("E16").Value = 1.002723585119218
$ws.Range("F16").Value = 0.9843644865386793
$ws.Range("G16").Value = 1.027139640605709
$ws.Range("H16").Value = 0.9445305089422378
$ws.Range("I16").Value = 1.011351016872909
$ws.Range("J16").Value = 0.9849234737043119
$ws.Range("K16").Value = 1.005664761589806
$ws.Range("L16").Value = 0.9635070909453787
$ws.Range("M16").Value = 1.027139640605709
$ws.Range("N16").Value = 0.9736270470307277
$ws.Range("O16").Value = 0.989689555301461
$ws.Range("P16").Value = 0.9905255705397811
